$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = 'LG 24GL600F-B '

$data = @(
    @('RUS', 'Официальный магазин LG', 491.56, '14 449 р.'),
    @('RUS', 'Goods.ru', 523.57, '15 390 р.'),
    @('RUS', 'Kns.ru', 546.78, '16 072 р.'),
    @('RUS', 'Svyaznoy.ru', 575.28, '16 910 р.'),
    @('RUS', 'CompYou.ru', 503.84, '14 810 р.'),
    @('RUS', 'Eldorado.ru', 523.57, '15 390 р.'),
    @('RUS', 'Elektro-park.ru', 522.55, '15 360 р.'),
    @('RUS', 'Citilink.ru', 503.84, '14 810 р.'),
    @('RUS', '123.ru', 533.1, '15 670 р.'),
    @('RUS', 'Империя Техно', 504.52, '14 830 р.'),
    @('RUS', 'М.видео', 523.57, '15 390 р.'),
    @('RUS', 'Pcplanet.ru', 557.63, '16 391 р.'),
    @('RUS', 'Topcomputer.ru', 475.6, '13 980 р.'),
    @('UA', 'Itbox.ua', 430.91, '4 777 грн.'),
    @('UA', 'LuxPRO.ua', 430.91, '4 777 грн.'),
    @('UA', 'Repka.ua', 431.27, '4 781 грн.'),
    @('UA', 'A-techno.com.ua', 430.91, '4 777 грн.'),
    @('UA', 'Homebt.com.ua', 451.03, '5 000 грн.'),
    @('UA', 'Isklad.tk', 430.91, '4 777 грн.'),
    @('UA', 'Brain.com.ua', 430.91, '4 777 грн.'),
    @('UA', 'Цифра', 425.41, '4 716 грн.'),
    @('UA', 'FOXTROT.UA', 430.91, '4 777 грн.'),
    @('UA', 'Епіцентр', 430.91, '4 777 грн.'),
    @('UA', 'Zhuk.ua', 430.91, '4 777 грн.'),
    @('UA', 'Stylus.ua', 425.41, '4 716 грн.'),
    @('UA', 'АЛЛО', 430.91, '4 777 грн.'),
    @('UA', 'Denika.ua', 430.91, '4 777 грн.'),
    @('UA', 'Palladium.ua', 430.91, '4 777 грн.'),
    @('UA', 'Hbox.com.ua', 466.09, '5 167 грн.'),
    @('UA', 'Rozetka.ua', 430.91, '4 777 грн.'),
    @('BLR', 'SOCKET.BY', 555.0, ''),
    @('BLR', 'TTN.by', 561.87, ''),
    @('BLR', 'KST.by', 555.0, ''),
    @('BLR', 'VIPCOMP.BY', 555.0, ''),
    @('BLR', 'Sli.by', 555.0, ''),
    @('BLR', 'Техник-Сервис', 576.42, ''),
    @('BLR', 'izliv.by', 555.0, ''),
    @('BLR', 'Медиа Маркет групп ООО', 582.37, ''),
    @('BLR', 'ЧУП "СВКомп"', 610.0, ''),
    @('BLR', 'ITMarket.by', 559.32, ''),
    @('BLR', 'Ньютон', 555.0, ''),
    @('BLR', 'Itsmart.by', 586.12, ''),
    @('BLR', 'AMD.by', 576.42, ''),
    @('BLR', 'imarket.by', 622.05, ''),
    @('BLR', 'ЭЛЕКТРОСИЛА', 622.0, ''),
    @('BLR', 'ViP MARKET', 606.39, ''),
    @('BLR', 'itx.by', 646.11, ''),
    @('BLR', 'BITS.By', 606.0, ''),
    @('BLR', 'ВИКО-ТЕХНО', 606.0, ''),
    @('BLR', 'Мультиком.', 556.68, ''),
    @('BLR', '5 элемент', 699.0, ''),
    @('BLR', 'bigi', 555.0, ''),
    @('BLR', '24shop.by', 555.0, ''),
    @('BLR', '21vek.by', 606.0, ''),
    @('BLR', 'LevelUP', 606.0, '')
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$ws.Range("A57:D60").Clear()
